# Add 5 new tax-rate columns (inserted between the existing columns) to the
# payoffs table, expanding the grid from B:G to B:L.
#
# Layout before:  B(=0) C(=1) D(=2) E(=3) F(=4) G(=5)
# Layout after:   B(=0) C(=1) D(=2) E(=3) F(=4) G(=5) H(=6) I(=7) J(=8) K(=9) L(=10)
# The original C,D,E,F,G values move to D,F,H,J,L respectively, and brand
# new values are written into C,E,G,I,K.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: extend the bold/bordered/centered style of the existing
# header cells (e.g. G1) onto the five new header cells, then fill values. ---
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:L1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("H1").Value = 6
$ws.Range("I1").Value = 7
$ws.Range("J1").Value = 8
$ws.Range("K1").Value = 9
$ws.Range("L1").Value = 10

# --- Data rows 2-7 ---
# The original values in C:G move out to D,F,H,J,L, and new values are
# inserted at C,E,G,I,K.
$data = @{
    2 = @{ C = 59.528;  D = 71.88;  E = 83.89700000000001; F = 98.053;  G = 114.98;  H = 115.604; I = 100.304; J = 85.004;  K = 69.70399999999999; L = 54.404 }
    3 = @{ C = 98.307;  D = 116.82; E = 132.38;             F = 147.463; G = 162.699; H = 156.404; I = 130.904; J = 105.404; K = 79.904;             L = 54.404 }
    4 = @{ C = 162.938; D = 191.72; E = 213.187;            F = 229.813; G = 242.23;  H = 224.404; I = 181.904; J = 139.404; K = 96.904;             L = 54.404 }
    5 = @{ C = 259.885; D = 304.07; E = 334.396;            F = 353.338; G = 361.527; H = 326.404; I = 258.404; J = 190.404; K = 122.404;            L = 54.404 }
    6 = @{ C = 518.41;  D = 603.67; E = 657.621;            F = 682.7380000000001; G = 679.652; H = 598.404; I = 462.404; J = 326.404; K = 190.404;  L = 54.404 }
    7 = @{ C = 809.251; D = 940.72; E = 1021.249;           F = 1053.313; G = 1037.543; H = 904.404; I = 691.904; J = 479.404; K = 266.904;         L = 54.404 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}
